# Weekly data refresh: a new record (week) is inserted at the top of the
# data block (row 639), pushing all existing records for this
# market/category down by one row. The former last row (730) therefore
# ends up duplicated into the new last row (731).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 639; Excel shifts rows 639:730 down to 640:731,
# carrying their values/styles with them (including the date column's
# style, and growing the sheet's used range / dimension to R731).
$ws.Rows.Item(639).Insert()

# Columns that are constant across this whole data block (market id,
# market name, region, codreg, category id, category, variety, unit,
# origin, kg/units flag, classification) — copy them from the row right
# below (the template for this series) so the new record matches the
# rest of the block.
$ws.Range("A639").Value2  = $ws.Range("A640").Value2
$ws.Range("B639").Value2  = $ws.Range("B640").Value2
$ws.Range("C639").Value2  = $ws.Range("C640").Value2
$ws.Range("E639").Value2  = $ws.Range("E640").Value2
$ws.Range("F639").Value2  = $ws.Range("F640").Value2
$ws.Range("G639").Value2  = $ws.Range("G640").Value2
$ws.Range("H639").Value2  = $ws.Range("H640").Value2
$ws.Range("N639").Value2  = $ws.Range("N640").Value2
$ws.Range("O639").Value2  = $ws.Range("O640").Value2
$ws.Range("Q639").Value2  = $ws.Range("Q640").Value2
$ws.Range("R639").Value2  = $ws.Range("R640").Value2

# New weekly record's own data.
$ws.Range("D639").Value2 = 45131
$ws.Range("I639").Value2 = "Primera"
$ws.Range("J639").Value2 = 650
$ws.Range("K639").Value2 = 900
$ws.Range("L639").Value2 = 1000
$ws.Range("M639").Value2 = 954
$ws.Range("P639").Value2 = 318
